$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure price/volume columns keep their original text representation
# (values like "42.814.85" or "  -0.48%  " are literal strings, not numbers/dates)
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "42.784.76"
$ws.Range("E2").Value = "  -0.67%  "
$ws.Range("D3").Value = "2.215.84"
$ws.Range("E3").Value = "  -1.16%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "256.08"
$ws.Range("E5").Value = "  +4.49%  "
$ws.Range("D6").Value = "0.617"
$ws.Range("E6").Value = "  +0.81%  "
$ws.Range("D7").Value = "77.59"
$ws.Range("E7").Value = "  +3.25%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("D9").Value = "0.593"
$ws.Range("E9").Value = "  -1.66%  "
$ws.Range("D10").Value = "42.19"
$ws.Range("E10").Value = "  +3.03%  "
$ws.Range("D11").Value = "0.0907"
$ws.Range("E11").Value = "  -2.42%  "
$ws.Range("D12").Value = "6.98"
$ws.Range("E12").Value = "  +0.98%  "
$ws.Range("E13").Value = "  +1.06%  "
$ws.Range("D14").Value = "2.548.72"
$ws.Range("E14").Value = "  -1.16%  "
$ws.Range("D15").Value = "14.44"
$ws.Range("E15").Value = "  -0.87%  "
$ws.Range("D16").Value = "2.213.29"
$ws.Range("E16").Value = "  -1.01%  "
$ws.Range("D17").Value = "0.779"
$ws.Range("E17").Value = "  -1.40%  "
$ws.Range("D18").Value = "42.773.33"
$ws.Range("E18").Value = "  -0.44%  "
$ws.Range("E19").Value = "  -1.85%  "
$ws.Range("D20").Value = "71.04"
$ws.Range("E20").Value = "  -0.13%  "
$ws.Range("E21").Value = "  -0.36%  "
$ws.Range("D22").Value = "2.30"
$ws.Range("E22").Value = "  +4.14%  "
$ws.Range("D23").Value = "230.03"
$ws.Range("E23").Value = "  +0.28%  "
$ws.Range("D24").Value = "9.29"
$ws.Range("E24").Value = "  -5.17%  "
$ws.Range("E25").Value = "  -0.18%  "
$ws.Range("E26").Value = "  +8.70%  "
$ws.Range("D27").Value = "10.77"
$ws.Range("E27").Value = "  -0.38%  "
$ws.Range("E28").Value = "  -2.46%  "
$ws.Range("B29").Value = "Toncoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D29").Value = "2.22"
$ws.Range("E29").Value = "  +3.88%  "
$ws.Range("B30").Value = "PancakeSwap"
$ws.Range("C30").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D30").Value = "2.19"
$ws.Range("E30").Value = "  -2.01%  "
$ws.Range("D31").Value = "173.27"
$ws.Range("E31").Value = "  +1.14%  "
$ws.Range("D32").Value = "20.43"
$ws.Range("E32").Value = "  +1.00%  "
$ws.Range("D33").Value = "0.0870"
$ws.Range("E33").Value = "  +8.43%  "
$ws.Range("E34").Value = "  -1.03%  "
$ws.Range("D35").Value = "0.121"
$ws.Range("E35").Value = "  -0.41%  "
$ws.Range("D36").Value = "0.0355"
$ws.Range("E36").Value = "  +7.76%  "
$ws.Range("E37").Value = "  -1.55%  "
$ws.Range("D38").Value = "4.33"
$ws.Range("E38").Value = "  -3.22%  "
$ws.Range("D39").Value = "13.22"
$ws.Range("E39").Value = "  +1.79%  "
$ws.Range("E40").Value = "  +17.46%  "
$ws.Range("E41").Value = "  -0.35%  "
$ws.Range("D42").Value = "0.200"
$ws.Range("E42").Value = "  -2.26%  "
$ws.Range("D43").Value = "5.32"
$ws.Range("E43").Value = "  -1.49%  "
$ws.Range("D44").Value = "60.50"
$ws.Range("E44").Value = "  +2.10%  "
$ws.Range("E45").Value = "  +0.32%  "
$ws.Range("D46").Value = "103.61"
$ws.Range("E46").Value = "  -1.00%  "
$ws.Range("D47").Value = "8.40"
$ws.Range("E47").Value = "  -3.32%  "
$ws.Range("E48").Value = "  -1.28%  "
$ws.Range("E49").Value = "  +0.84%  "
$ws.Range("E50").Value = "  -1.10%  "
$ws.Range("B51").Value = "Stacks"
$ws.Range("C51").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D51").Value = "1.46"
$ws.Range("E51").Value = "  +18.81%  "
